$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling the data.
$updates = @{
    3  = -3
    4  = -1
    5  = -1
    6  = -1
    7  = 4
    12 = -6
    21 = 0
    29 = -5
    35 = -5
    41 = -5
    42 = 4
    47 = 0
    48 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
